$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column widths (closest achievable values under this host's rounding)
$ws.Columns.Item(1).ColumnWidth = 16.333333333333332
$ws.Columns.Item(2).ColumnWidth = 14
$ws.Columns.Item(3).ColumnWidth = 11

# Header row
$ws.Range("A1").Value = "Match_factor_floor"
$ws.Range("B1").Value = "Reverse_match_factor_floor"
$ws.Range("C1").Value = "LRI_diff_floor"
$ws.Range("D1").Value = "mean_pct_volfound"
$ws.Range("E1").Value = "median_pct_found"

# Column A - Match_factor_floor ladder (600, 625, then +25 each row)
$ws.Range("A2").Value = 600
$ws.Range("A3").Value = 625
$ws.Range("A4").Formula = "=A3+25"
$ws.Range("A5:A11").Formula = "=A4+25"
$ws.Range("A10").Formula = "=A9+25"

# Column B - Reverse_match_factor_floor (constant 100, carried via formula)
$ws.Range("B2").Value = 100
$ws.Range("B3").Formula = "=B2"
$ws.Range("B4:B11").Formula = "=B3"

# Column C - LRI_diff_floor (constant 10, carried via formula)
$ws.Range("C2").Value = 10
$ws.Range("C3").Formula = "=C2"
$ws.Range("C4:C11").Formula = "=C3"

# Column D - mean_pct_volfound (sparse values)
$ws.Range("D2").Value = 65.5
$ws.Range("D3").Value = 64.72
$ws.Range("D4").Value = 63.7
$ws.Range("D5").Value = 62.3
$ws.Range("D6").Value = 60.6
$ws.Range("D8").Value = 55.9
$ws.Range("D10").Value = 48.959

# Column E - median_pct_found (sparse values)
$ws.Range("E2").Value = 70.9
$ws.Range("E3").Value = 70.2
$ws.Range("E4").Value = 69
$ws.Range("E5").Value = 67.7
$ws.Range("E6").Value = 65.4
$ws.Range("E8").Value = 60.5
$ws.Range("E10").Value = 52.99

# Row 12 - new literal row (not continuing the A:C formula chain)
$ws.Range("A12").Value = 750
$ws.Range("B12").Value = 100
$ws.Range("C12").Value = 5
$ws.Range("D12").Value = 48.6
$ws.Range("E12").Value = 52.54

# Final selection lands one row below the data, matching the authored file
$ws.Range("E13").Select()
